$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1804281345565749
$ws.Range("C2").Value = 0.5688073394495413
$ws.Range("J2").Value = 0.01529051987767584
$ws.Range("P2").Value = 0.1284403669724771
$ws.Range("S2").Value = 0.1070336391437309
$ws.Range("B3").Value = 0.01047120418848168
$ws.Range("C3").Value = 0.01047120418848168
$ws.Range("J3").Value = 0.05759162303664921
$ws.Range("P3").Value = 0.7225130890052356
$ws.Range("S3").Value = 0.1989528795811518
$ws.Range("J4").Value = 0.09090909090909091
$ws.Range("P4").Value = 0.7045454545454546
$ws.Range("S4").Value = 0.2045454545454546
$ws.Range("B6").Value = 0.04508196721311476
$ws.Range("D6").Value = 0.02459016393442623
$ws.Range("F6").Value = 0.06557377049180328
$ws.Range("J6").Value = 0.290983606557377
$ws.Range("O6").Value = 0.02049180327868852
$ws.Range("Q6").Value = 0.1352459016393443
$ws.Range("R6").Value = 0.06967213114754098
$ws.Range("S6").Value = 0.3483606557377049
$ws.Range("B7").Value = 0.1004784688995215
$ws.Range("D7").Value = 0.01435406698564593
$ws.Range("E7").Value = 0.004784688995215311
$ws.Range("F7").Value = 0.08133971291866028
$ws.Range("J7").Value = 0.1244019138755981
$ws.Range("O7").Value = 0.03827751196172249
$ws.Range("Q7").Value = 0.1866028708133971
$ws.Range("R7").Value = 0.07177033492822966
$ws.Range("S7").Value = 0.3779904306220095
$ws.Range("B8").Value = 0.1056466302367942
$ws.Range("D8").Value = 0.01457194899817851
$ws.Range("F8").Value = 0.0546448087431694
$ws.Range("J8").Value = 0.122040072859745
$ws.Range("O8").Value = 0.03642987249544627
$ws.Range("Q8").Value = 0.1930783242258652
$ws.Range("R8").Value = 0.07832422586520947
$ws.Range("S8").Value = 0.395264116575592
$ws.Range("B9").Value = 0.09236947791164658
$ws.Range("D9").Value = 0.008032128514056224
$ws.Range("F9").Value = 0.05220883534136546
$ws.Range("J9").Value = 0.07228915662650602
$ws.Range("O9").Value = 0.02008032128514056
$ws.Range("Q9").Value = 0.2008032128514056
$ws.Range("R9").Value = 0.08433734939759036
$ws.Range("S9").Value = 0.4698795180722892
$ws.Range("B10").Value = 0.105643994211288
$ws.Range("D10").Value = 0.01881331403762663
$ws.Range("F10").Value = 0.06512301013024602
$ws.Range("J10").Value = 0.1418234442836469
$ws.Range("O10").Value = 0.03111432706222865
$ws.Range("Q10").Value = 0.2076700434153401
$ws.Range("R10").Value = 0.06367583212735166
$ws.Range("S10").Value = 0.3661360347322721
$ws.Range("G11").Value = 0.157556270096463
$ws.Range("J11").Value = 0.09646302250803858
$ws.Range("K11").Value = 0.2090032154340836
$ws.Range("L11").Value = 0.5176848874598071
$ws.Range("S11").Value = 0.01929260450160772
$ws.Range("G12").Value = 0.7337278106508875
$ws.Range("J12").Value = 0.1538461538461539
$ws.Range("K12").Value = 0.005917159763313609
$ws.Range("L12").Value = 0.04142011834319527
$ws.Range("S12").Value = 0.0650887573964497
$ws.Range("G13").Value = 0.78
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.02
$ws.Range("F15").Value = 0.0119047619047619
$ws.Range("H15").Value = 0.1388888888888889
$ws.Range("I15").Value = 0.05952380952380952
$ws.Range("J15").Value = 0.2738095238095238
$ws.Range("K15").Value = 0.06746031746031746
$ws.Range("M15").Value = 0.01984126984126984
$ws.Range("N15").Value = 0.003968253968253968
$ws.Range("O15").Value = 0.03174603174603174
$ws.Range("S15").Value = 0.3928571428571428
$ws.Range("F16").Value = 0.01507537688442211
$ws.Range("H16").Value = 0.2412060301507538
$ws.Range("I16").Value = 0.09547738693467336
$ws.Range("J16").Value = 0.3869346733668342
$ws.Range("K16").Value = 0.1005025125628141
$ws.Range("M16").Value = 0.03015075376884422
$ws.Range("O16").Value = 0.03015075376884422
$ws.Range("S16").Value = 0.1005025125628141
$ws.Range("F17").Value = 0.01565557729941291
$ws.Range("H17").Value = 0.1800391389432485
$ws.Range("I17").Value = 0.09001956947162426
$ws.Range("J17").Value = 0.3894324853228963
$ws.Range("K17").Value = 0.09784735812133072
$ws.Range("M17").Value = 0.01956947162426614
$ws.Range("N17").Value = 0.003913894324853229
$ws.Range("O17").Value = 0.05870841487279843
$ws.Range("S17").Value = 0.1448140900195695
$ws.Range("F18").Value = 0.04891304347826087
$ws.Range("H18").Value = 0.1684782608695652
$ws.Range("I18").Value = 0.1141304347826087
$ws.Range("J18").Value = 0.4130434782608696
$ws.Range("K18").Value = 0.07608695652173914
$ws.Range("M18").Value = 0.02717391304347826
$ws.Range("O18").Value = 0.03804347826086957
$ws.Range("S18").Value = 0.1141304347826087
$ws.Range("F19").Value = 0.01995870612525809
$ws.Range("H19").Value = 0.2381280110116999
$ws.Range("I19").Value = 0.09910529938059189
$ws.Range("J19").Value = 0.3688919476944253
$ws.Range("K19").Value = 0.09291121816930489
$ws.Range("M19").Value = 0.01789401238816242
$ws.Range("O19").Value = 0.05918788713007571
$ws.Range("S19").Value = 0.1039229181004818
